# The deck's slide master currently uses the "Integral" theme colour
# scheme (ppt/theme/theme1.xml). The commit swaps it for the stock
# "Office Theme" colour scheme (the set of 12 theme colours that used
# to sit unused in ppt/theme/theme2.xml, wired only to the notes
# master). Re-point every theme colour slot on the master to the
# Office Theme RGB values via ThemeColorScheme.Colors(i).RGB, which is
# the supported way to edit a deck's theme colours through the object
# model.

$p  = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

# Office Theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# RGB() packs as 0x00BBGGRR, matching ColorFormat.RGB's native order.
$cs.Colors(1).RGB  = 0x000000   # Dark 1    - 000000
$cs.Colors(2).RGB  = 0xFFFFFF   # Light 1   - FFFFFF
$cs.Colors(3).RGB  = 0x6A5444   # Dark 2    - 44546A
$cs.Colors(4).RGB  = 0xE6E6E7   # Light 2   - E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B   # Accent 1  - 5B9BD5
$cs.Colors(6).RGB  = 0x317DED   # Accent 2  - ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5   # Accent 3  - A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF   # Accent 4  - FFC000
$cs.Colors(9).RGB  = 0xC47244   # Accent 5  - 4472C4
$cs.Colors(10).RGB = 0x47AD70   # Accent 6  - 70AD47
$cs.Colors(11).RGB = 0xC16305   # Hyperlink - 0563C1
$cs.Colors(12).RGB = 0x724F95   # Followed  - 954F72
